{"js": "// Locate the existing \"_GoBack\" bookmark paragraph (the second paragraph\n// in the body). We insert a new misspelled-word paragraph directly before\n// it, and a new empty paragraph directly after it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst bookmarkParagraph = paragraphs.items[1];\n\n// --- New paragraph inserted BEFORE the bookmark paragraph -----------------\n// It must contain two separate runs (\"W\" + \"tfyughoijop\") bracketed by\n// w:proofErr spell-check markers, matching the target OOXML exactly.\n// insertParagraph(\"\", \"Before\") gives us a placeholder paragraph that we\n// then replace in-place via insertOoxml (Flat OPC wrapped WordprocessingML)\n// so the run/proofErr structure is preserved verbatim (plain insertText\n// would merge same-format runs into a single <w:r>).\nconst misspelledParagraph = bookmarkParagraph.insertParagraph(\"\", \"Before\");\n\nconst misspelledFlatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>W</w:t></w:r>\n            <w:r><w:t>tfyughoijop</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nmisspelledParagraph.insertOoxml(misspelledFlatOpc, \"Replace\");\n\n// --- New empty paragraph inserted AFTER the bookmark paragraph ------------\nconst trailingParagraph = bookmarkParagraph.insertParagraph(\"\", \"After\");\n\nconst emptyFlatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p/>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntrailingParagraph.insertOoxml(emptyFlatOpc, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Locate the existing \"_GoBack\" bookmark paragraph (the second paragraph\n# in the document body). A new misspelled-word paragraph is inserted\n# directly before it, and a new empty paragraph is inserted directly\n# after it.\n$d = $word.ActiveDocument\n\n$bookmarkParagraph = $d.Paragraphs.Item(2).Range\n\n# --- New paragraph inserted BEFORE the bookmark paragraph ------------------\n# It must contain two separate runs (\"W\" + \"tfyughoijop\") bracketed by\n# w:proofErr spell-check markers, matching the target OOXML exactly.\n# InsertParagraphBefore() gives us a placeholder paragraph that we then\n# replace in-place via Range.InsertXML so the run/proofErr structure is\n# preserved verbatim (plain Range.Text would merge into a single run with\n# no proofErr markers).\n$bookmarkParagraph.InsertParagraphBefore()\n$misspelledRange = $d.Paragraphs.Item(2).Range\n$misspelledXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:proofErr w:type=\"spellStart\"/><w:r><w:t>W</w:t></w:r><w:r><w:t>tfyughoijop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n$misspelledRange.InsertXML($misspelledXml)\n\n# --- New empty paragraph inserted AFTER the bookmark paragraph ------------\n$bookmarkParagraph2 = $d.Paragraphs.Item(3).Range\n$bookmarkParagraph2.InsertParagraphAfter()\n$trailingRange = $d.Paragraphs.Item(4).Range\n$emptyXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"/>'\n$trailingRange.InsertXML($emptyXml)\n"}
